$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "28.383.31"
Set-TextValue "E2" "  +0.34%  "
Set-TextValue "D3" "1.869.75"
Set-TextValue "E3" "  +0.03%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "330.27"
Set-TextValue "E5" "  -2.83%  "
Set-TextValue "D6" "1.000"
Set-TextValue "E6" "  -0.07%  "
Set-TextValue "D7" "0.4611"
Set-TextValue "E7" "  -1.98%  "
Set-TextValue "D8" "0.4012"
Set-TextValue "E8" "  +2.21%  "
Set-TextValue "D9" "47.74"
Set-TextValue "E9" "  +0.91%  "
Set-TextValue "D10" "0.07853"
Set-TextValue "E10" "  -1.87%  "
Set-TextValue "D11" "0.9854"
Set-TextValue "E11" "  -2.01%  "
Set-TextValue "D12" "21.31"
Set-TextValue "E12" "  -2.56%  "
Set-TextValue "D13" "1.869.71"
Set-TextValue "E13" "  -0.10%  "
Set-TextValue "D14" "5.852"
Set-TextValue "E14" "  -2.50%  "
Set-TextValue "D15" "6.991"
Set-TextValue "E15" "  -3.97%  "
Set-TextValue "E16" "  +0.00%  "
Set-TextValue "D17" "88.12"
Set-TextValue "E17" "  -3.39%  "
Set-TextValue "D18" "0.06537"
Set-TextValue "E18" "  -1.04%  "
Set-TextValue "E19" "  -2.37%  "
Set-TextValue "D20" "17.22"
Set-TextValue "E20" "  -2.33%  "
Set-TextValue "D21" "0.9998"
Set-TextValue "E21" "  -0.19%  "
Set-TextValue "D22" "28.368.37"
Set-TextValue "E22" "  +0.29%  "
Set-TextValue "D23" "5.343"
Set-TextValue "E23" "  -1.93%  "
Set-TextValue "D24" "10.86"
Set-TextValue "E24" "  -1.85%  "
Set-TextValue "D25" "2.248"
Set-TextValue "E25" "  -1.81%  "
Set-TextValue "D26" "2.090.53"
Set-TextValue "E26" "  -0.16%  "
Set-TextValue "D27" "157.84"
Set-TextValue "E27" "  -1.30%  "
Set-TextValue "D28" "19.38"
Set-TextValue "E28" "  -2.46%  "
Set-TextValue "D29" "2.060"
Set-TextValue "E29" "  -4.05%  "
Set-TextValue "D30" "5.287"
Set-TextValue "E30" "  -4.01%  "
Set-TextValue "D31" "117.48"
Set-TextValue "E31" "  -2.31%  "
Set-TextValue "D32" "0.9587"
Set-TextValue "E32" "  -1.98%  "
Set-TextValue "D33" "0.09340"
Set-TextValue "E33" "  -1.88%  "
Set-TextValue "E34" "  -0.16%  "
Set-TextValue "E35" "  +0.76%  "
Set-TextValue "D36" "5.239"
Set-TextValue "E36" "  -2.14%  "
Set-TextValue "D37" "0.06030"
Set-TextValue "E37" "  -1.07%  "
Set-TextValue "D38" "0.02199"
Set-TextValue "E38" "  -3.09%  "
Set-TextValue "D39" "8.289"
Set-TextValue "E39" "  -1.63%  "
Set-TextValue "D40" "1.159"
Set-TextValue "E40" "  -1.72%  "
Set-TextValue "D41" "0.9999"
Set-TextValue "E41" "  -0.11%  "
Set-TextValue "D42" "0.5759"
Set-TextValue "E42" "  -3.63%  "
Set-TextValue "D43" "0.1807"
Set-TextValue "E43" "  -4.00%  "
Set-TextValue "D44" "10.03"
Set-TextValue "E44" "  -3.52%  "
Set-TextValue "D45" "1.244"
Set-TextValue "E45" "  -3.19%  "
Set-TextValue "D46" "2.299"
Set-TextValue "E46" "  +13.53%  "
Set-TextValue "D47" "0.5430"
Set-TextValue "E47" "  -3.31%  "
Set-TextValue "D48" "11.87"
Set-TextValue "E48" "  -2.31%  "
Set-TextValue "D49" "0.07134"
Set-TextValue "E49" "  +3.29%  "
Set-TextValue "D50" "1.886"
Set-TextValue "E50" "  -4.18%  "
Set-TextValue "D51" "111.02"
Set-TextValue "E51" "  -0.27%  "
